# Replace the test staff member name on the "StaffMember" sheet and leave
# that sheet as the active/selected one (matches the saved workbook state).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("StaffMember")
$ws.Activate()
$ws.Range("A2").Value = "Michael Goldstein"
$ws.Range("B7").Select()
